# Update values in column C for the specified rows on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    12  = -10.9696
    27  = -12.2018
    32  = -12.4286
    36  = -11.84240000000001
    38  = -12.332
    46  = -14.68509999999999
    54  = -12.90980000000001
    55  = -13.54269999999999
    56  = -12.7806
    67  = -12.0081
    69  = -11.7908
    72  = -11.82480000000001
    83  = -13.57
    86  = -14.25389999999999
    91  = -12.3829
    93  = -10.4811
    99  = -11.99930000000001
    104 = -12.34420000000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
